$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the planning-time timestamps (column A) ---
$ws.Range("A2").Value = 45229.7291666667
$ws.Range("A3").Value = 45229.7708333333
$ws.Range("A4").Value = 45229.75

# --- Row 2: new product folder/description + new image list ---
$ws.Range("C2").Value = "/Users/jishuliu/Desktop/hkshop/data/20231030_superdelivery/【日本直送】 ササガワ　ＩＴ事業部（すべてのジャンル）  荷札シール　取扱注意  行李標籤貼紙 小心輕放  一色入"
$ws.Range("D2").Value = "https://shopage.s3.amazonaws.com/media/f854/615273998674_40870190220846178232.jpg,https://shopage.s3.amazonaws.com/media/f854/615273998674_22676963380593410709.jpg"

# --- Row 3: same description, new image list ---
$ws.Range("C3").Value = "/Users/jishuliu/Desktop/hkshop/data/20231030_superdelivery/【日本直送】 ササガワ　ＩＴ事業部（すべてのジャンル）  荷札シール　取扱注意  行李標籤貼紙 小心輕放  一色入"
$ws.Range("D3").Value = "https://shopage.s3.amazonaws.com/media/f854/615273998674_03156904136975129854.jpg,https://shopage.s3.amazonaws.com/media/f854/615273998674_65403258431629315721.jpg"

# --- Row 4: same description, new image list (E4 already holds an empty string; leave as-is) ---
$ws.Range("C4").Value = "/Users/jishuliu/Desktop/hkshop/data/20231030_superdelivery/【日本直送】 ササガワ　ＩＴ事業部（すべてのジャンル）  荷札シール　取扱注意  行李標籤貼紙 小心輕放  一色入"
$ws.Range("D4").Value = "https://shopage.s3.amazonaws.com/media/f854/615273998674_46386897831229038267.jpg,https://shopage.s3.amazonaws.com/media/f854/615273998674_35055401077700880486.jpg"

# --- Column width adjustments to fit the new content ---
# (ColumnWidth is stored internally in whole-pixel increments of 1/7 character
#  units, so the chosen values are the closest achievable to the target widths
#  of 27.3942307692308 / 174 / 244.5 once re-serialized.)
$ws.Columns.Item(1).ColumnWidth = 26.714285714285715
$ws.Columns.Item(3).ColumnWidth = 173.28571428571428
$ws.Columns.Item(4).ColumnWidth = 243.71428571428572

# --- Restore the last active selection ---
$ws.Range("B7").Select()
